# Auto-update predictions and index for 2025-10-29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (rows 8 and 9) - their fixtures have rolled off the tracked list
$ws.Rows("8:9").Delete()

# Row 3: PAOK Thessaloniki - Volos NPS finished 3:0, prediction correct
$ws.Range("A3").Value = "PAOK Thessaloniki ✓ - Volos NPS: 3:0"
$ws.Range("B3").Value = "PAOK Thessaloniki"
$ws.Range("C3").Value = 72
$ws.Range("D3").Value = 86
$ws.Range("F3").Value = 1.29
$ws.Range("G3").Value = "✓"

# Row 4: Real Madrid - FC Barcelona finished 2:1, prediction correct
$ws.Range("A4").Value = "Real Madrid ✓ - FC Barcelona: 2:1"
$ws.Range("B4").Value = "Real Madrid"
$ws.Range("C4").Value = 70
$ws.Range("D4").Value = 49
$ws.Range("F4").Value = 2.05
$ws.Range("G4").Value = "✓"

# Row 5: CD Tondela - Sporting CP, still pending, updated confidence figures
$ws.Range("A5").Value = "CD Tondela - Sporting CP : 17:00"
$ws.Range("B5").Value = "Sporting CP"
$ws.Range("C5").Value = 69
$ws.Range("D5").Value = 100
$ws.Range("F5").Value = 1.28

# Row 6: FK Bodø/Glimt - Molde FK finished 4:1, prediction correct
$ws.Range("A6").Value = "FK Bodø/Glimt ✓ - Molde FK: 4:1"
$ws.Range("B6").Value = "FK Bodø/Glimt"
$ws.Range("C6").Value = 64
$ws.Range("D6").Value = 84
$ws.Range("F6").Value = 1.29
$ws.Range("G6").Value = "✓"

# Row 7: CE Casa de Portugal - FC Santa Coloma B, still pending, updated confidence figures
$ws.Range("A7").Value = "CE Casa de Portugal  - FC Santa Coloma B: 13:00"
$ws.Range("B7").Value = "CE Casa de Portugal"
$ws.Range("C7").Value = 55
$ws.Range("D7").Value = 22
# Oddspedia_Confidence no longer available for this fixture -> blank it out
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = 13

Write-Host "applied edits"
